# UndoRedoSequenceDiagram.pptx
# DeveloperGuide.adoc: Update implementation for undo/redo to OrderBook
#
# The deck's single slide hosts a UML sequence diagram that talks about an
# "AddressBook" model; it needs to become "OrderBook" instead. Four shapes
# carry the old wording:
#   - a lifeline header reading ":Address" / "BookParser" (two paragraphs)
#   - a call label "undoAddressBook()"
#   - a lifeline header ":VersionedAddressBook"
#   - a call label "resetData(ReadOnlyAddressBook)"
#
# Rather than hard-coding shape indices, walk every shape on the slide and
# patch only the exact substring that changed (via TextRange.Characters),
# so the untouched runs/paragraphs keep their original formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Replace-Substring($shapes, $needle, $replacement) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $hasTf = $false
        try { $hasTf = [bool]$sh.HasTextFrame } catch { $hasTf = $false }
        if (-not $hasTf) { continue }

        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $pos = $full.IndexOf($needle)
        if ($pos -ge 0) {
            $tr.Characters($pos + 1, $needle.Length).Text = $replacement
        }
    }
}

# Order matters: replace the longer/more specific names first so that a
# shorter needle ("AddressBook") can't accidentally match inside a shape
# that should be matched by a longer one first. Needles are kept narrow
# (e.g. "AddressBook", not "undoAddressBook") so a replacement never spans
# two differently-formatted runs and collapses their formatting together.
Replace-Substring $s.Shapes "VersionedAddressBook" "VersionedOrderBook"
Replace-Substring $s.Shapes "ReadOnlyAddressBook" "ReadOnlyOrderBook"
Replace-Substring $s.Shapes "AddressBook" "OrderBook"
Replace-Substring $s.Shapes ":Address" ":Order"
